# update scripts wuth new tpm
# Refresh the Hbegf-Cd44 LR-pair sheet with recomputed TPM-derived statistics
# (ligand/receptor expression values and downstream specificity/weight metrics).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.94498
$ws.Range("H2").Value = 17.83494
$ws.Range("I2").Value = 0.4679240463447598
$ws.Range("J2").Value = 0.4679240463447597
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 44.51376303756
$ws.Range("R2").Value = 400.62386733804
$ws.Range("S2").Value = 0.06477917940073338
$ws.Range("T2").Value = 0.06477917940073337
$ws.Range("G3").Value = 5.94498
$ws.Range("H3").Value = 17.83494
$ws.Range("I3").Value = 0.4679240463447598
$ws.Range("J3").Value = 0.4679240463447597
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 190.23595947144
$ws.Range("R3").Value = 1712.12363524296
$ws.Range("S3").Value = 0.2768431268475961
$ws.Range("T3").Value = 0.276843126847596
$ws.Range("G4").Value = 5.94498
$ws.Range("H4").Value = 17.83494
$ws.Range("I4").Value = 0.4679240463447598
$ws.Range("J4").Value = 0.4679240463447597
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 86.7897028319
$ws.Range("R4").Value = 781.1073254871
$ws.Range("S4").Value = 0.1263017400964303
$ws.Range("T4").Value = 0.1263017400964303
$ws.Range("H5").Value = 8.352077
$ws.Range("I5").Value = 0.219128164447035
$ws.Range("J5").Value = 0.219128164447035
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 20.845731830298
$ws.Range("R5").Value = 187.611586472682
$ws.Range("S5").Value = 0.03033599744948618
$ws.Range("T5").Value = 0.03033599744948618
$ws.Range("H6").Value = 8.352077
$ws.Range("I6").Value = 0.219128164447035
$ws.Range("J6").Value = 0.219128164447035
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("Q6").Value = 89.08722887065198
$ws.Range("R6").Value = 801.785059835868
$ws.Range("S6").Value = 0.1296452419997987
$ws.Range("T6").Value = 0.1296452419997987
$ws.Range("H7").Value = 8.352077
$ws.Range("I7").Value = 0.219128164447035
$ws.Range("J7").Value = 0.219128164447035
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 40.64349422308944
$ws.Range("R7").Value = 365.7914480078049
$ws.Range("S7").Value = 0.05914692499775009
$ws.Range("T7").Value = 0.05914692499775009
$ws.Range("G8").Value = 3.976005
$ws.Range("H8").Value = 11.928015
$ws.Range("I8").Value = 0.3129477892082053
$ws.Range("J8").Value = 0.3129477892082053
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 29.77082251011
$ws.Range("R8").Value = 267.93740259099
$ws.Range("S8").Value = 0.04332434107317652
$ws.Range("T8").Value = 0.04332434107317652
$ws.Range("G9").Value = 3.976005
$ws.Range("H9").Value = 11.928015
$ws.Range("I9").Value = 0.3129477892082053
$ws.Range("J9").Value = 0.3129477892082053
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 127.22988572514
$ws.Range("R9").Value = 1145.06897152626
$ws.Range("S9").Value = 0.1851527938801604
$ws.Range("T9").Value = 0.1851527938801604
$ws.Range("G10").Value = 3.976005
$ws.Range("H10").Value = 11.928015
$ws.Range("I10").Value = 0.3129477892082053
$ws.Range("J10").Value = 0.3129477892082053
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 58.044987940775
$ws.Range("R10").Value = 522.404891466975
$ws.Range("S10").Value = 0.08447065425486835
$ws.Range("T10").Value = 0.08447065425486835
